$p = $ppt.ActivePresentation
$s = $p.Slides.Item(37)

# --- Change 1: update trailing punctuation in the End User Testing paragraph ---
$shp2 = $s.Shapes.Item(2)
$tr = $shp2.TextFrame.TextRange
$paraRange = $tr.Paragraphs(7,1)

# Step 1: write full target text but with junk markers on both ends so the host
# treats this as an entirely new run (keeps our own right Unicode punctuation).
$step1 = "ZZWhen a new round is wanting to  be played the code previously stated that “press <enter> to play again but if any thing except ‘x’ will play another round. This has now been changed to show that any key can be entered to play again:ZZ"
$paraRange.Text = $step1

# Step 2: strip the junk markers with another whole-both-ends-differ replace so it
# again becomes a single run (instead of splitting prefix/suffix runs).
$tr2 = $shp2.TextFrame.TextRange
$paraRange2 = $tr2.Paragraphs(7,1)
$step2 = "When a new round is wanting to  be played the code previously stated that “press <enter> to play again but if any thing except ‘x’ will play another round. This has now been changed to show that any key can be entered to play again:"
$paraRange2.Text = $step2

# --- Change 2: move/resize the Rectangle 22 highlight shape ---
$shp5 = $s.Shapes.Item(5)
$shp5.Left = 314.30767822265625
$shp5.Top = 434.14459228515625
$shp5.Width = 111.60035705566406
$shp5.Height = 20.932401657104492
